# Fruta / hortaliza, semanal
# Insert a new weekly record row before row 96 (shifts existing rows 96-99 down to 97-100)
# and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96; existing row 96 (and below) shift down to 97.
$ws.Rows.Item(96).Insert()

# New row 96 shares the same categorical/descriptive values as the surrounding rows
# (Mercado ID, Mercado, Region, Codreg, Tipo, Producto ID, Producto, Categoria ID,
#  Categoria, Variedad, Calidad, Unidad de comercializacion, Origen, Kg/unidad),
# only the date, volume and prices differ for this week.
$ws.Cells.Item(96, 1).Value = 10
$ws.Cells.Item(96, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(96, 3).Value = "La Araucanía"
$ws.Cells.Item(96, 4).Value = 45021
$ws.Cells.Item(96, 4).NumberFormat = $ws.Cells.Item(97, 4).NumberFormat
$ws.Cells.Item(96, 5).Value = 9
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100108
$ws.Cells.Item(96, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(96, 9).Value = 100108004
$ws.Cells.Item(96, 10).Value = "Papaya"
$ws.Cells.Item(96, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 80
$ws.Cells.Item(96, 14).Value = 27000
$ws.Cells.Item(96, 15).Value = 27000
$ws.Cells.Item(96, 16).Value = 27000
$ws.Cells.Item(96, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(96, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(96, 19).Value = 2700
$ws.Cells.Item(96, 20).Value = 10
